# Week06/CssSelectorsAndProperties.pptx - "first round of changes I implemented"
#
# Slide 6 ("Ruleset Example") - Content Placeholder 2:
#   - the glitch.com link text/URL gains a "/remix" path segment
#   - both runs in that paragraph (the link run and the trailing space run)
#   get a light-gray (C0C0C0) highlight applied

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)

# Update just the first run's text so the run split / hyperlink (rId3) on
# the first run is preserved exactly as-is.
$tr = $shp.TextFrame.TextRange
$r1 = $tr.Runs(1)
$r1.Text = "https://glitch.com/edit/#!/remix/rulesetexample"

# Apply a silver (C0C0C0) highlight to every run in the text box (the link
# run and the following space run) via the DrawingML-level text range.
$tf2 = $shp.TextFrame2
$tr2 = $tf2.TextRange
$tr2.Font.Highlight.RGB = 12632256
